$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay text even though it looks numeric
# (Excel/COM would otherwise silently coerce a bare numeric-looking string
# into a number). Cells whose new text is unambiguously non-numeric (two
# dots, trailing letters, special glyphs, surrounding spaces, etc.) are
# set with a plain .Value assignment instead, so we don't touch their
# number format unnecessarily.
function Set-NumericLookingText($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.492.40"
$ws.Range("E2").Value = "  -2.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.801.03"
$ws.Range("E3").Value = "  -2.40%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.41%  "

# Row 5 - BNB
Set-NumericLookingText "D5" "229.07"
$ws.Range("E5").Value = "  -1.60%  "

# Row 6 - XRP
Set-NumericLookingText "D6" "0.611"
$ws.Range("E6").Value = "  -1.28%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.43%  "

# Row 8 - Solana
Set-NumericLookingText "D8" "39.33"
$ws.Range("E8").Value = "  -11.78%  "

# Row 9 - Cardano
Set-NumericLookingText "D9" "0.321"
$ws.Range("E9").Value = "  +2.92%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.01%  "

# Row 11 - TRON
Set-NumericLookingText "D11" "0.0988"
$ws.Range("E11").Value = "  -2.20%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.059.83"
$ws.Range("E12").Value = "  -2.40%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  -2.28%  "

# Row 14 - now Polygon (was WrappedEther)
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-NumericLookingText "D14" "0.658"
$ws.Range("E14").Value = "  -2.69%  "

# Row 15 - now WrappedEther (was Polygon)
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.790.08"
$ws.Range("E15").Value = "  -2.97%  "

# Row 16 - Polkadot
Set-NumericLookingText "D16" "4.55"
$ws.Range("E16").Value = "  -3.95%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "34.332.51"
$ws.Range("E17").Value = "  -3.35%  "

# Row 18 - Litecoin
Set-NumericLookingText "D18" "68.98"
$ws.Range("E18").Value = "  -2.30%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0779"
$ws.Range("E19").Value = "  -2.99%  "

# Row 20 - BitcoinCash
Set-NumericLookingText "D20" "239.10"
$ws.Range("E20").Value = "  -2.28%  "

# Row 21 - Avalanche
Set-NumericLookingText "D21" "11.76"
$ws.Range("E21").Value = "  -2.88%  "

# Row 22 - Uniswap
Set-NumericLookingText "D22" "4.68"
$ws.Range("E22").Value = "  +0.87%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.36%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -1.71%  "

# Row 25 - Monero
Set-NumericLookingText "D25" "173.04"
$ws.Range("E25").Value = "  +1.04%  "

# Row 26 - Cosmos
Set-NumericLookingText "D26" "7.73"
$ws.Range("E26").Value = "  -3.76%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -3.91%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.44%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  -5.59%  "

# Row 30 - BinanceUSD
$ws.Range("E30").Value = "  +0.58%  "

# Row 31 - Filecoin
Set-NumericLookingText "D31" "3.99"
$ws.Range("E31").Value = "  +1.05%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -2.26%  "

# Row 33 - InternetComputer(DFINITY)
Set-NumericLookingText "D33" "3.90"
$ws.Range("E33").Value = "  -5.72%  "

# Row 34 - TrustWalletToken
Set-NumericLookingText "D34" "1.23"
$ws.Range("E34").Value = "  +7.73%  "

# Row 35 - LidoDAOToken
Set-NumericLookingText "D35" "1.79"
$ws.Range("E35").Value = "  -3.16%  "

# Row 36 - ImmutableX
Set-NumericLookingText "D36" "0.692"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37 - Aave
Set-NumericLookingText "D37" "90.39"
$ws.Range("E37").Value = "  -5.49%  "

# Row 38 - WEMIXToken
$ws.Range("E38").Value = "  +5.00%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.323.69"
$ws.Range("E39").Value = "  -2.02%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -3.01%  "

# Row 41 - ARBITRUM
Set-NumericLookingText "D41" "0.959"
$ws.Range("E41").Value = "  -6.01%  "

# Row 42 - now InjectiveProtocol (was HuobiToken)
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-NumericLookingText "D42" "14.13"
$ws.Range("E42").Value = "  -7.83%  "

# Row 43 - now HuobiToken (was InjectiveProtocol)
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-NumericLookingText "D43" "2.38"
$ws.Range("E43").Value = "  -3.03%  "

# Row 44 - RenderToken
Set-NumericLookingText "D44" "2.21"
$ws.Range("E44").Value = "  -9.55%  "

# Row 45 - MXToken
$ws.Range("E45").Value = "  -3.71%  "

# Row 46 - FraxShare
$ws.Range("E46").Value = "  -1.97%  "

# Row 47 - Kaspa
$ws.Range("E47").Value = "  -1.40%  "

# Row 48 - RocketPoolETH
$ws.Range("D48").Value = "1.983.39"
$ws.Range("E48").Value = "  -1.57%  "

# Row 50 - Cronos
Set-NumericLookingText "D50" "0.0658"
$ws.Range("E50").Value = "  +3.54%  "

# Row 51 - Quant
Set-NumericLookingText "D51" "97.56"
$ws.Range("E51").Value = "  -4.94%  "
